$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1: "obs", formatted like the other header cells (e.g. A1)
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "obs"

# New data cells K2:K6: "Ape", formatted like the other plain data cells (e.g. J2)
$ws.Range("J2").Copy()
$ws.Range("K2:K6").PasteSpecial(-4122)
$ws.Range("K2:K6").Value = "Ape"

$excel.CutCopyMode = 0

# Match the final selection recorded in the workbook
$ws.Range("K2:K6").Select()
